$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 13 (original data rows that are removed)
$ws.Range("A4:A13").EntireRow.Delete()

# Add new header "change" in D1, copying formatting from C1 (bold header style)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "change"

# Update row 2 values
$ws.Range("B2").Value = "2024-06-15 20:49:16.736937"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 0

# Update row 3 values
$ws.Range("B3").Value = "2024-06-15 20:49:46.269381"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = -2
